# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values for each game row (row 2 = most recent game, ... row 27 = oldest)
$kValues = @{
    2  = 7
    3  = 7
    4  = 2
    5  = 4
    6  = 4
    7  = 6
    8  = 5
    9  = 10
    10 = 7
    11 = 12
    12 = 5
    13 = 3
    14 = 6
    15 = 5
    16 = 3
    17 = 6
    18 = 5
    19 = 8
    20 = 7
    21 = 7
    22 = 8
    23 = 3
    24 = 4
    25 = 5
    26 = 4
    27 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
